# Delete row 138 ("「美しい字は目で見る音楽」" entry) from Sheet1.
# This shifts all subsequent rows up by one, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Rows.Item(138).Delete()
